{"js": "// The edit does three things to \"Planificacion del proyecto.docx\":\n//   1. Relocates the `_GoBack` bookmark from the empty trailing paragraph\n//      (after the table) to a spot inside the \"Sexto entregable -\n//      Incremento de software 2 (sprint 2)\" table cell, right after the\n//      word \"Incremento\" (this is where the author's cursor/selection was\n//      when the document was last saved).\n//   2. Changes \"Modulo compras y reportes\" to \"Modulo ventas y reportes\"\n//      in the table.\n//   3. As a consequence of (1), the trailing paragraph after the table no\n//      longer carries the bookmark.\n\n// --- 1) Relocate the \"_GoBack\" bookmark -----------------------------------\n// Remove it from wherever it currently lives (the empty paragraph that\n// follows the table).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Find the split point: right after \"Sexto entregable - Incremento\" and\n// before \" de software 2 (sprint 2)\".\nconst titleResults = context.document.body.search(\n  \"Sexto entregable - Incremento\",\n  { matchCase: true, matchWholeWord: false }\n);\ntitleResults.load(\"text\");\nawait context.sync();\n\nconst splitPoint = titleResults.items[0].getRange(\"End\");\nsplitPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2) \"Modulo compras y reportes\" -> \"Modulo ventas y reportes\" ---------\n// Scope the search to the specific table cell so we don't touch the\n// unrelated \"Modulo de compras\" row elsewhere in the table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst targetRow = table.rows.items[27];\ntargetRow.cells.load(\"items\");\nawait context.sync();\n\nconst targetCell = targetRow.cells.items[0];\nconst wordResults = targetCell.body.search(\"compras\", { matchCase: true });\nwordResults.load(\"text\");\nawait context.sync();\n\nwordResults.items[0].insertText(\"ventas\", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The edit does three things to \"Planificacion del proyecto.docx\":\n#   1. Relocates the `_GoBack` bookmark from the empty trailing paragraph\n#      (after the table) to a spot inside the \"Sexto entregable -\n#      Incremento de software 2 (sprint 2)\" table cell, right after the\n#      word \"Incremento\" (this is where the author's cursor/selection was\n#      when the document was last saved).\n#   2. Changes \"Modulo compras y reportes\" to \"Modulo ventas y reportes\"\n#      in the table.\n#   3. As a consequence of (1), the trailing paragraph after the table no\n#      longer carries the bookmark.\n\n# --- 1) Relocate the \"_GoBack\" bookmark -----------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$titleRng = $d.Content\n$titleRng.Find.Execute(\"Sexto entregable - Incremento\") | Out-Null\n$titleRng.Collapse(0)  # wdCollapseEnd: collapse to the end of the found text\n$d.Bookmarks.Add(\"_GoBack\", $titleRng) | Out-Null\n\n# --- 2) \"Modulo compras y reportes\" -> \"Modulo ventas y reportes\" ---------\n# Locate the unique phrase first (there is another, unrelated \"compras\" in\n# the document, in a \"Modulo de compras\" row) and then swap just the word\n# \"compras\" for \"ventas\" within it.\n$bodyRng = $d.Content\n$bodyRng.Find.Execute(\"Modulo compras y reportes\") | Out-Null\n\n$prefixLen = \"Modulo \".Length\n$wordLen = \"compras\".Length\n$wordStart = $bodyRng.Start + $prefixLen\n$wordEnd = $wordStart + $wordLen\n\n$wordRng = $d.Range($wordStart, $wordEnd)\n$wordRng.Text = \"ventas\"\n"}
